$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 858.6667
$ws.Range("I40").Value = 812.8570999999999
$ws.Range("J40").Value = 1500
$ws.Range("K40").Value = 812.8570999999999
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = -637.8570999999999
$ws.Range("N40").Value = -1850
$ws.Range("H82").Value = 2973.6667
$ws.Range("I82").Value = 1868.4
$ws.Range("J82").Value = 8500
$ws.Range("K82").Value = 5605.200000000001
$ws.Range("L82").Value = 25500
$ws.Range("M82").Value = -5199.200000000001
$ws.Range("N82").Value = -26312
$ws.Range("H85").Value = 2973.6667
$ws.Range("I85").Value = 1868.4
$ws.Range("J85").Value = 8500
$ws.Range("K85").Value = 5605.200000000001
$ws.Range("L85").Value = 25500
$ws.Range("M85").Value = -4201.200000000001
$ws.Range("N85").Value = -28308
$ws.Range("H127").Value = 758.9286
$ws.Range("I127").Value = 755.7692
$ws.Range("J127").Value = 800
$ws.Range("K127").Value = 2267.3076
$ws.Range("L127").Value = 2400
$ws.Range("M127").Value = 2692.6924
$ws.Range("N127").Value = -12320
$ws.Range("H129").Value = 1241.1923
$ws.Range("I129").Value = 397
$ws.Range("J129").Value = 1274.96
$ws.Range("K129").Value = 1191
$ws.Range("L129").Value = 3824.88
$ws.Range("M129").Value = 3809
$ws.Range("N129").Value = -13824.88
$ws.Range("H137").Value = 1720.65
$ws.Range("I137").Value = 1221.2727
$ws.Range("J137").Value = 2331
$ws.Range("K137").Value = 3663.8181
$ws.Range("L137").Value = 6993
$ws.Range("M137").Value = -1113.8181
$ws.Range("N137").Value = -12093
$ws.Range("H138").Value = 4785.234
$ws.Range("I138").Value = 4018.7334
$ws.Range("J138").Value = 5144.5312
$ws.Range("K138").Value = 12056.2002
$ws.Range("L138").Value = 15433.5936
$ws.Range("M138").Value = -6916.200199999999
$ws.Range("N138").Value = -25713.5936
$ws.Range("H141").Value = 4094.4211
$ws.Range("I141").Value = 4180.875
$ws.Range("J141").Value = 3633.3333
$ws.Range("K141").Value = 12542.625
$ws.Range("L141").Value = 10899.9999
$ws.Range("M141").Value = -7362.625
$ws.Range("N141").Value = -21259.9999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 772935.2
$ws.Range("I32").Value = 9604.603999999999
$ws.Range("J32").Value = 10887065
$ws.Range("K32").Value = 9604.603999999999
$ws.Range("L32").Value = 10887065
$ws.Range("M32").Value = -9317.603999999999
$ws.Range("N32").Value = -10887639
$ws.Range("H61").Value = 3462.4666
$ws.Range("J61").Value = 2514
$ws.Range("L61").Value = 2514
$ws.Range("N61").Value = -2938
$ws.Range("H74").Value = 1066.3334
$ws.Range("I74").Value = 1194.6666
$ws.Range("J74").Value = 873.8333
$ws.Range("K74").Value = 1194.6666
$ws.Range("L74").Value = 873.8333
$ws.Range("M74").Value = -320.6666
$ws.Range("N74").Value = -2621.8333
$ws.Range("H77").Value = 1066.3334
$ws.Range("I77").Value = 1194.6666
$ws.Range("J77").Value = 873.8333
$ws.Range("K77").Value = 5973.333000000001
$ws.Range("L77").Value = 4369.1665
$ws.Range("M77").Value = -1605.333000000001
$ws.Range("N77").Value = -13105.1665
$ws.Range("H122").Value = 3238.8044
$ws.Range("I122").Value = 3258.225
$ws.Range("K122").Value = 9774.674999999999
$ws.Range("M122").Value = -7324.674999999999
$ws.Range("H136").Value = 3462.4666
$ws.Range("J136").Value = 2514
$ws.Range("L136").Value = 7542
$ws.Range("N136").Value = -12642

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7591.8
$ws.Range("I86").Value = 1993
$ws.Range("J86").Value = 15990
$ws.Range("K86").Value = 1993
$ws.Range("L86").Value = 15990
$ws.Range("M86").Value = -870
$ws.Range("N86").Value = -18236
$ws.Range("H89").Value = 7591.8
$ws.Range("I89").Value = 1993
$ws.Range("J89").Value = 15990
$ws.Range("K89").Value = 9965
$ws.Range("L89").Value = 79950
$ws.Range("M89").Value = -4349
$ws.Range("N89").Value = -91182
$ws.Range("H94").Value = 942.06665
$ws.Range("I94").Value = 942.06665
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 942.06665
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -491.06665
$ws.Range("H139").Value = 20000
$ws.Range("J139").Value = 20000
$ws.Range("L139").Value = 20000
$ws.Range("N139").Value = -30280

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10364.25
$ws.Range("I31").Value = 3685.25
$ws.Range("K31").Value = 3685.25
$ws.Range("M31").Value = -3390.25
$ws.Range("H34").Value = 10364.25
$ws.Range("I34").Value = 3685.25
$ws.Range("K34").Value = 3685.25
$ws.Range("M34").Value = -3483.25
$ws.Range("H58").Value = 1398.174
$ws.Range("I58").Value = 810.73334
$ws.Range("J58").Value = 2499.625
$ws.Range("K58").Value = 810.73334
$ws.Range("L58").Value = 2499.625
$ws.Range("M58").Value = -607.73334
$ws.Range("N58").Value = -2905.625
$ws.Range("H134").Value = 861.8570999999999
$ws.Range("I134").Value = 835.8
$ws.Range("J134").Value = 927
$ws.Range("K134").Value = 2507.4
$ws.Range("L134").Value = 2781
$ws.Range("M134").Value = 27.60000000000036
$ws.Range("N134").Value = -7851
$ws.Range("H136").Value = 1398.174
$ws.Range("I136").Value = 810.73334
$ws.Range("J136").Value = 2499.625
$ws.Range("K136").Value = 2432.20002
$ws.Range("L136").Value = 7498.875
$ws.Range("M136").Value = 117.7999799999998
$ws.Range("N136").Value = -12598.875

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 880.14
$ws.Range("I113").Value = 494.16666
$ws.Range("J113").Value = 932.7727
$ws.Range("K113").Value = 1482.49998
$ws.Range("L113").Value = 2798.3181
$ws.Range("M113").Value = 687.5000199999999
$ws.Range("N113").Value = -7138.3181
$ws.Range("H122").Value = 475.89743
$ws.Range("J122").Value = 1300.8
$ws.Range("L122").Value = 11707.2
$ws.Range("N122").Value = -16607.2
$ws.Range("H131").Value = 11905900
$ws.Range("J131").Value = 13158826
$ws.Range("L131").Value = 39476478
$ws.Range("N131").Value = -39486558

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 818
$ws.Range("I2").Value = 30
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 30
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = 82
$ws.Range("N2").Value = -2224
$ws.Range("H139").Value = 79550
$ws.Range("J139").Value = 79550
$ws.Range("L139").Value = 79550
$ws.Range("N139").Value = -89830

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1105.4667
$ws.Range("I136").Value = 631.8333
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 1895.4999
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = 654.5001
$ws.Range("N136").Value = -14100
$ws.Range("H139").Value = 64650
$ws.Range("J139").Value = 64650
$ws.Range("L139").Value = 64650
$ws.Range("N139").Value = -74930
